# Auto-generated Excel COM-interop script
# Refreshes Universalis market-board price snapshots (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ) plus their dependent Leve
# price / profit columns in each class table, matching the scheduled-runner
# commit's re-pull of live prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1870.3914
$ws.Range("I15").Value = 1870.3914
$ws.Range("K15").Value = 5611.174199999999
$ws.Range("M15").Value = -5442.174199999999
$ws.Range("H40").Value = 2501.8096
$ws.Range("I40").Value = 1887
$ws.Range("J40").Value = 2646.4707
$ws.Range("K40").Value = 1887
$ws.Range("L40").Value = 2646.4707
$ws.Range("M40").Value = -1712
$ws.Range("N40").Value = -2996.4707
$ws.Range("H116").Value = 16206.613
$ws.Range("J116").Value = 15696.296
$ws.Range("L116").Value = 15696.296
$ws.Range("N116").Value = -22580.296
$ws.Range("H132").Value = 3741
$ws.Range("I132").Value = 3594.5117
$ws.Range("K132").Value = 10783.5351
$ws.Range("M132").Value = -8253.535100000001
$ws.Range("H135").Value = 1216.16
$ws.Range("I135").Value = 1216.16
$ws.Range("K135").Value = 10945.44
$ws.Range("M135").Value = -8410.440000000001
$ws.Range("H137").Value = 1794.2
$ws.Range("I137").Value = 797
$ws.Range("J137").Value = 2610.0908
$ws.Range("K137").Value = 2391
$ws.Range("L137").Value = 7830.2724
$ws.Range("M137").Value = 159
$ws.Range("N137").Value = -12930.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6521.293
$ws.Range("I61").Value = 5351.1797
$ws.Range("J61").Value = 8923.105
$ws.Range("K61").Value = 5351.1797
$ws.Range("L61").Value = 8923.105
$ws.Range("M61").Value = -5139.1797
$ws.Range("N61").Value = -9347.105
$ws.Range("H74").Value = 4261.3335
$ws.Range("I74").Value = 3945.9092
$ws.Range("K74").Value = 3945.9092
$ws.Range("M74").Value = -3071.9092
$ws.Range("H77").Value = 4261.3335
$ws.Range("I77").Value = 3945.9092
$ws.Range("K77").Value = 19729.546
$ws.Range("M77").Value = -15361.546
$ws.Range("H88").Value = 4849.0713
$ws.Range("I88").Value = 761.9
$ws.Range("J88").Value = 7119.722
$ws.Range("K88").Value = 761.9
$ws.Range("L88").Value = 7119.722
$ws.Range("M88").Value = -355.9
$ws.Range("N88").Value = -7931.722
$ws.Range("H91").Value = 4849.0713
$ws.Range("I91").Value = 761.9
$ws.Range("J91").Value = 7119.722
$ws.Range("K91").Value = 761.9
$ws.Range("L91").Value = 7119.722
$ws.Range("M91").Value = 642.1
$ws.Range("N91").Value = -9927.722
$ws.Range("H136").Value = 6521.293
$ws.Range("I136").Value = 5351.1797
$ws.Range("J136").Value = 8923.105
$ws.Range("K136").Value = 16053.5391
$ws.Range("L136").Value = 26769.315
$ws.Range("M136").Value = -13503.5391
$ws.Range("N136").Value = -31869.315

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1334433.5
$ws.Range("I22").Value = 1082.091
$ws.Range("J22").Value = 5001150
$ws.Range("K22").Value = 1082.091
$ws.Range("L22").Value = 5001150
$ws.Range("M22").Value = -909.0909999999999
$ws.Range("N22").Value = -5001496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 874.75
$ws.Range("I3").Value = 750
$ws.Range("K3").Value = 750
$ws.Range("M3").Value = -637
$ws.Range("H31").Value = 2850.9092
$ws.Range("I31").Value = 1712
$ws.Range("J31").Value = 3346.087
$ws.Range("K31").Value = 1712
$ws.Range("L31").Value = 3346.087
$ws.Range("M31").Value = -1417
$ws.Range("N31").Value = -3936.087
$ws.Range("H34").Value = 2850.9092
$ws.Range("I34").Value = 1712
$ws.Range("J34").Value = 3346.087
$ws.Range("K34").Value = 1712
$ws.Range("L34").Value = 3346.087
$ws.Range("M34").Value = -1510
$ws.Range("N34").Value = -3750.087
$ws.Range("H58").Value = 3417.5356
$ws.Range("I58").Value = 2459.05
$ws.Range("K58").Value = 2459.05
$ws.Range("M58").Value = -2256.05
$ws.Range("H136").Value = 3417.5356
$ws.Range("I136").Value = 2459.05
$ws.Range("K136").Value = 7377.150000000001
$ws.Range("M136").Value = -4827.150000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3548.3333
$ws.Range("I3").Value = 2278
$ws.Range("J3").Value = 9900
$ws.Range("K3").Value = 6834
$ws.Range("L3").Value = 29700
$ws.Range("M3").Value = -6722
$ws.Range("N3").Value = -29924
$ws.Range("H7").Value = 599.75
$ws.Range("J7").Value = 749.8333
$ws.Range("L7").Value = 2249.4999
$ws.Range("N7").Value = -2473.4999
$ws.Range("H62").Value = 7957.5
$ws.Range("J62").Value = 8299
$ws.Range("L62").Value = 24897
$ws.Range("N62").Value = -26269
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H65").Value = 7957.5
$ws.Range("J65").Value = 8299
$ws.Range("L65").Value = 74691
$ws.Range("N65").Value = -81555
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H68").Value = 3531
$ws.Range("I68").Value = 1124.75
$ws.Range("K68").Value = 3374.25
$ws.Range("M68").Value = -2563.25
$ws.Range("H71").Value = 3531
$ws.Range("I71").Value = 1124.75
$ws.Range("K71").Value = 10122.75
$ws.Range("M71").Value = -6066.75
$ws.Range("H130").Value = 12028.833
$ws.Range("I130").Value = 5670
$ws.Range("K130").Value = 17010
$ws.Range("M130").Value = -11990
$ws.Range("H131").Value = 3646.068
$ws.Range("I131").Value = 2198.3333
$ws.Range("K131").Value = 6594.999899999999
$ws.Range("M131").Value = -1554.999899999999
$ws.Range("H133").Value = 22400.305
$ws.Range("I133").Value = 8151.25
$ws.Range("K133").Value = 24453.75
$ws.Range("M133").Value = -19393.75
$ws.Range("H134").Value = 9118.799999999999
$ws.Range("I134").Value = 3453.6667
$ws.Range("J134").Value = 17616.5
$ws.Range("K134").Value = 10361.0001
$ws.Range("L134").Value = 52849.5
$ws.Range("M134").Value = -5291.000100000001
$ws.Range("N134").Value = -62989.5
$ws.Range("H137").Value = 7118.7334
$ws.Range("I137").Value = 3091.75
$ws.Range("J137").Value = 11721
$ws.Range("K137").Value = 9275.25
$ws.Range("L137").Value = 35163
$ws.Range("M137").Value = -4175.25
$ws.Range("N137").Value = -45363
$ws.Range("H138").Value = 20711.924
$ws.Range("I138").Value = 3295.3333
$ws.Range("J138").Value = 29932.47
$ws.Range("K138").Value = 9885.999899999999
$ws.Range("L138").Value = 89797.41
$ws.Range("M138").Value = -4745.999899999999
$ws.Range("N138").Value = -100077.41
$ws.Range("H139").Value = 8982.700000000001
$ws.Range("I139").Value = 1685.4
$ws.Range("J139").Value = 16280
$ws.Range("K139").Value = 5056.200000000001
$ws.Range("L139").Value = 48840
$ws.Range("M139").Value = 83.79999999999927
$ws.Range("N139").Value = -59120
$ws.Range("H140").Value = 9620858
$ws.Range("I140").Value = 14707754
$ws.Range("J140").Value = 12277.556
$ws.Range("K140").Value = 44123262
$ws.Range("L140").Value = 36832.66800000001
$ws.Range("M140").Value = -44118082
$ws.Range("N140").Value = -47192.66800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6532.28
$ws.Range("I70").Value = 6476.8945
$ws.Range("K70").Value = 6476.8945
$ws.Range("M70").Value = -6206.8945
$ws.Range("H73").Value = 6532.28
$ws.Range("I73").Value = 6476.8945
$ws.Range("K73").Value = 6476.8945
$ws.Range("M73").Value = -5540.8945

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1619.2
$ws.Range("I22").Value = 598
$ws.Range("J22").Value = 1874.5
$ws.Range("K22").Value = 598
$ws.Range("L22").Value = 1874.5
$ws.Range("M22").Value = -303
$ws.Range("N22").Value = -2464.5
$ws.Range("H27").Value = 1619.2
$ws.Range("I27").Value = 598
$ws.Range("J27").Value = 1874.5
$ws.Range("K27").Value = 598
$ws.Range("L27").Value = 1874.5
$ws.Range("M27").Value = -491
$ws.Range("N27").Value = -2088.5
$ws.Range("H46").Value = 2143.111
$ws.Range("I46").Value = 1066
$ws.Range("J46").Value = 2450.8572
$ws.Range("K46").Value = 1066
$ws.Range("L46").Value = 2450.8572
$ws.Range("M46").Value = -878
$ws.Range("N46").Value = -2826.8572
$ws.Range("H76").Value = 36249.75
$ws.Range("J76").Value = 36249.75
$ws.Range("L76").Value = 36249.75
$ws.Range("N76").Value = -36925.75
$ws.Range("H79").Value = 36249.75
$ws.Range("J79").Value = 36249.75
$ws.Range("L79").Value = 36249.75
$ws.Range("N79").Value = -38589.75
$ws.Range("H136").Value = 2198.7778
$ws.Range("I136").Value = 1825.4286
$ws.Range("J136").Value = 3505.5
$ws.Range("K136").Value = 5476.2858
$ws.Range("L136").Value = 10516.5
$ws.Range("M136").Value = -2926.2858
$ws.Range("N136").Value = -15616.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4252.35
$ws.Range("I132").Value = 4221.9414
$ws.Range("K132").Value = 12665.8242
$ws.Range("M132").Value = -10135.8242
